$d = $word.ActiveDocument

# Locate the target text " <fr>en noyau</fr>" within the document.
$found = $d.Content.Duplicate
$ok = $found.Find.Execute(" <fr>en noyau</fr>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($ok) {
    $s = $found.Start
    $e = $found.End

    # Sub-ranges for each piece of the split run:
    #   " "        -> unchanged formatting
    #   "<fr>"     -> Courier New, blue, size 9pt (sz 18 half-points)
    #   "en noyau" -> unchanged formatting
    #   "</fr>"    -> Courier New, blue, size 9pt (sz 18 half-points)
    $rOpenTag = $d.Range($s + 1, $s + 5)
    $rCloseTag = $d.Range($s + 13, $e)

    $rOpenTag.Font.Name = "Courier New"
    $rOpenTag.Font.Color = 255
    $rOpenTag.Font.Size = 9

    $rCloseTag.Font.Name = "Courier New"
    $rCloseTag.Font.Color = 255
    $rCloseTag.Font.Size = 9
}
